$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1. "target" category list (column A): drop the stray "text" entry ---
# that used to sit between "step" and "web", shifting the remaining
# categories (web, webalert, webcookie, ws, ws.async, xml) up one row.
$ws.Range("A25").Delete("xlShiftUp")

# --- 2. Function list columns (Y..AE): the old standalone "text" column
# (Y) goes away, so web/webalert/webcookie/ws/ws.async/xml each shift one
# column to the left (Z->Y, AA->Z, AB->AA, AC->AB, AD->AC, AE->AD).
$ws.Columns("Y").Delete()

# --- 3. "json" function list (column M): add the new storeKeys(...)
# function in alphabetical order, right before storeValue(...). Shift the
# two trailing entries down a row by hand (scoped to column M only) and
# drop the new entry into the gap this opens up.
$ws.Range("M18").Value2 = $ws.Range("M17").Value2
$ws.Range("M17").Value2 = $ws.Range("M16").Value2
$ws.Range("M16").Value2 = "storeKeys(json,jsonpath,var)"

# --- 4. Keep the workbook-level named ranges in sync with the new
# locations/sizes of the lists above (the automation layer does not
# recompute these automatically the way interactive Excel would).
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
